$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.395131349563599
$ws.Range("B1").Value = 3.33491039276123
$ws.Range("C1").Value = 5.148038864135742
$ws.Range("D1").Value = 7.214332580566406
$ws.Range("E1").Value = 3.952540159225464
